$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '98.688.43'
$ws.Range('E2').Value = '  +0.74%  '
$ws.Range('D3').Value = '3.452.59'
$ws.Range('E3').Value = '  +1.59%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '255.72'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.78%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '674.48'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.60%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.51'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.78%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.435'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.42%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.07'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.12%  '
$ws.Range('E10').Value = '  -0.02%  '
$ws.Range('D11').Value = '3.444.88'
$ws.Range('E11').Value = '  +1.47%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '47.09'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +13.70%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.212'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.17%  '
$ws.Range('D14').Value = '98.546.46'
$ws.Range('E14').Value = '  +0.92%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.21'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.22%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000261'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.80%  '
$ws.Range('D17').Value = '4.095.51'
$ws.Range('E17').Value = '  +1.73%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '9.16'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.76%  '
$ws.Range('D19').Value = '3.438.14'
$ws.Range('E19').Value = '  +1.21%  '
$ws.Range('E20').Value = '  +7.39%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.539'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.44%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.80'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +8.00%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '523.76'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.68%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.46'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.85%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000204'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.60%  '
$ws.Range('E26').Value = '  +6.13%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '98.47'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.64%  '
$ws.Range('E28').Value = '  +1.71%  '
$ws.Range('D29').Value = '3.623.89'
$ws.Range('E29').Value = '  +1.22%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.06'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +18.18%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '12.33'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +7.73%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.147'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.75%  '
$ws.Range('E33').Value = '  +0.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.191'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.30%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.14%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.578'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.86%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '29.86'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.80%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.21'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.18%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.54'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.42%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '538.88'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.93%  '
$ws.Range('E41').Value = '  +2.97%  '
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.881'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.90%  '
$ws.Range('B44').Value = 'ImmutableX'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.81'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.51%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0439'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.60%  '
$ws.Range('B46').Value = 'WhiteBITCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '24.44'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.09%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.75'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.63%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.77'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.02%  '
$ws.Range('E49').Value = '  -2.65%  '
$ws.Range('E50').Value = '  +7.72%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '55.78'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.59%  '
